$wb = $excel.ActiveWorkbook

# --- Rename "NL Stats-this session" to "PLO Stats-this session" ---
$wsPLO = $wb.Worksheets.Item("NL Stats-this session")
$wsPLO.Name = "PLO Stats-this session"

$wsCombined = $wb.Worksheets.Item("combined Stats-this session")

$pctFormat = "0.0%"
$dollarFormat = "`"$`"#,##0.00_-"

function Set-StatsRow {
    param($ws, $row, $data)
    foreach ($col in $data.Keys) {
        $cellRef = "$col$row"
        $kind = $data[$col][0]
        $val = $data[$col][1]
        $cell = $ws.Range($cellRef)
        if ($kind -eq "ds") {
            # Date-looking string: prefix with apostrophe so it is stored as
            # literal text instead of being parsed into a date serial number.
            $cell.Value = "'" + $val
        } else {
            $cell.Value = $val
        }
        if ($kind -eq "p") {
            $cell.NumberFormat = $pctFormat
        } elseif ($kind -eq "d") {
            $cell.NumberFormat = $dollarFormat
        }
    }
}

$row2 = @{
    "A" = @("s", "Fish")
    "B" = @("n", 20)
    "C" = @("n", 26.14)
    "D" = @("n", 6.14)
    "E" = @("n", 0)
    "F" = @("p", 0.711)
    "G" = @("p", 0.018)
    "H" = @("p", 0)
    "I" = @("p", 0.21)
    "J" = @("p", 0.303)
    "K" = @("p", 0.165)
    "L" = @("n", 0.37)
    "M" = @("n", 0)
    "N" = @("n", 2)
    "O" = @("d", 161.46)
    "P" = @("d", 79.63)
    "Q" = @("n", 218)
    "R" = @("p", 0.545)
    "T" = @("s", "Date")
}

$row3 = @{
    "A" = @("s", "Cedric")
    "B" = @("n", 20)
    "C" = @("n", 12.39)
    "D" = @("n", -7.61)
    "E" = @("n", 0)
    "F" = @("p", 0.586)
    "G" = @("p", 0)
    "H" = @("p", 0)
    "I" = @("p", 0.266)
    "J" = @("p", 0.276)
    "K" = @("p", 0.124)
    "L" = @("n", 0.63)
    "M" = @("n", 0)
    "N" = @("n", 0)
    "O" = @("d", 159.28)
    "P" = @("d", 54.6)
    "Q" = @("n", 210)
    "R" = @("p", 0.448)
    "T" = @("ds", "05/20/21")
}

$row4 = @{
    "A" = @("s", "Scott")
    "B" = @("n", 20)
    "C" = @("n", 55.27)
    "D" = @("n", 35.27)
    "E" = @("n", 0)
    "F" = @("p", 0.657)
    "G" = @("p", 0.343)
    "H" = @("p", 0.0178)
    "I" = @("p", 0.523)
    "J" = @("p", 0.16)
    "K" = @("p", 0.107)
    "L" = @("n", 3.02)
    "M" = @("n", 42)
    "N" = @("n", 75)
    "O" = @("d", 190.59)
    "P" = @("d", 96.95999999999999)
    "Q" = @("n", 169)
    "R" = @("p", 0.667)
}

$row5 = @{
    "A" = @("s", "Xavier")
    "B" = @("n", 50)
    "C" = @("n", 0)
    "D" = @("n", -50)
    "E" = @("n", 2)
    "F" = @("p", 0.847)
    "G" = @("p", 0.492)
    "H" = @("p", 0)
    "I" = @("p", 0.465)
    "J" = @("p", 0.458)
    "K" = @("p", 0.136)
    "L" = @("n", 1.18)
    "M" = @("n", 11)
    "N" = @("n", 25)
    "O" = @("d", 74.38)
    "P" = @("d", 27.1)
    "Q" = @("n", 59)
    "R" = @("p", 0.296)
}

foreach ($ws in @($wsCombined, $wsPLO)) {
    Set-StatsRow $ws 2 $row2
    Set-StatsRow $ws 3 $row3
    Set-StatsRow $ws 4 $row4
    Set-StatsRow $ws 5 $row5
}
# --- Update chart series source ranges ---
function Set-ChartSeriesFormula {
    param($ws, $chartObjIndex, $seriesIndex, $formula)
    $co = $ws.ChartObjects().Item($chartObjIndex)
    $ser = $co.Chart.SeriesCollection().Item($seriesIndex)
    $ser.Formula = $formula
}

# chart1.xml -> wsCombined ChartObject 1
Set-ChartSeriesFormula $wsCombined 1 1 "=SERIES('combined Stats-this session'!F1,'combined Stats-this session'!`$A`$2:`$A`$5,'combined Stats-this session'!`$F`$2:`$F`$5,1)"
Set-ChartSeriesFormula $wsCombined 1 2 "=SERIES('combined Stats-this session'!G1,'combined Stats-this session'!`$A`$2:`$A`$5,'combined Stats-this session'!`$G`$2:`$G`$5,2)"
Set-ChartSeriesFormula $wsCombined 1 3 "=SERIES('combined Stats-this session'!H1,'combined Stats-this session'!`$A`$2:`$A`$5,'combined Stats-this session'!`$H`$2:`$H`$5,3)"

# chart2.xml -> wsCombined ChartObject 2
Set-ChartSeriesFormula $wsCombined 2 1 "=SERIES('combined Stats-this session'!L1,'combined Stats-this session'!`$A`$2:`$A`$5,'combined Stats-this session'!`$L`$2:`$L`$5,1)"

# chart3.xml -> wsCombined ChartObject 3
Set-ChartSeriesFormula $wsCombined 3 1 "=SERIES('combined Stats-this session'!M1,'combined Stats-this session'!`$A`$2:`$A`$5,'combined Stats-this session'!`$M`$2:`$M`$5,1)"
Set-ChartSeriesFormula $wsCombined 3 2 "=SERIES('combined Stats-this session'!N1,'combined Stats-this session'!`$A`$2:`$A`$5,'combined Stats-this session'!`$N`$2:`$N`$5,2)"

# chart4.xml -> wsCombined ChartObject 4
Set-ChartSeriesFormula $wsCombined 4 1 "=SERIES('combined Stats-this session'!J1,'combined Stats-this session'!`$A`$2:`$A`$5,'combined Stats-this session'!`$J`$2:`$J`$5,1)"
Set-ChartSeriesFormula $wsCombined 4 2 "=SERIES('combined Stats-this session'!K1,'combined Stats-this session'!`$A`$2:`$A`$5,'combined Stats-this session'!`$K`$2:`$K`$5,2)"

# chart5.xml -> wsCombined ChartObject 5
Set-ChartSeriesFormula $wsCombined 5 1 "=SERIES('combined Stats-this session'!R1,'combined Stats-this session'!`$A`$2:`$A`$5,'combined Stats-this session'!`$R`$2:`$R`$5,1)"

# chart6.xml -> wsCombined ChartObject 6
Set-ChartSeriesFormula $wsCombined 6 1 "=SERIES('combined Stats-this session'!O1,'combined Stats-this session'!`$A`$2:`$A`$5,'combined Stats-this session'!`$O`$2:`$O`$5,1)"
Set-ChartSeriesFormula $wsCombined 6 2 "=SERIES('combined Stats-this session'!P1,'combined Stats-this session'!`$A`$2:`$A`$5,'combined Stats-this session'!`$P`$2:`$P`$5,2)"

# chart7.xml -> wsCombined ChartObject 7
Set-ChartSeriesFormula $wsCombined 7 1 "=SERIES('combined Stats-this session'!Q1,'combined Stats-this session'!`$A`$2:`$A`$5,'combined Stats-this session'!`$Q`$2:`$Q`$5,1)"

# chart8.xml -> wsPLO ChartObject 1
Set-ChartSeriesFormula $wsPLO 1 1 "=SERIES('PLO Stats-this session'!F1,'PLO Stats-this session'!`$A`$2:`$A`$5,'PLO Stats-this session'!`$F`$2:`$F`$5,1)"
Set-ChartSeriesFormula $wsPLO 1 2 "=SERIES('PLO Stats-this session'!G1,'PLO Stats-this session'!`$A`$2:`$A`$5,'PLO Stats-this session'!`$G`$2:`$G`$5,2)"
Set-ChartSeriesFormula $wsPLO 1 3 "=SERIES('PLO Stats-this session'!H1,'PLO Stats-this session'!`$A`$2:`$A`$5,'PLO Stats-this session'!`$H`$2:`$H`$5,3)"

# chart9.xml -> wsPLO ChartObject 2
Set-ChartSeriesFormula $wsPLO 2 1 "=SERIES('PLO Stats-this session'!L1,'PLO Stats-this session'!`$A`$2:`$A`$5,'PLO Stats-this session'!`$L`$2:`$L`$5,1)"

# chart10.xml -> wsPLO ChartObject 3
Set-ChartSeriesFormula $wsPLO 3 1 "=SERIES('PLO Stats-this session'!M1,'PLO Stats-this session'!`$A`$2:`$A`$5,'PLO Stats-this session'!`$M`$2:`$M`$5,1)"
Set-ChartSeriesFormula $wsPLO 3 2 "=SERIES('PLO Stats-this session'!N1,'PLO Stats-this session'!`$A`$2:`$A`$5,'PLO Stats-this session'!`$N`$2:`$N`$5,2)"

# chart11.xml -> wsPLO ChartObject 4
Set-ChartSeriesFormula $wsPLO 4 1 "=SERIES('PLO Stats-this session'!J1,'PLO Stats-this session'!`$A`$2:`$A`$5,'PLO Stats-this session'!`$J`$2:`$J`$5,1)"
Set-ChartSeriesFormula $wsPLO 4 2 "=SERIES('PLO Stats-this session'!K1,'PLO Stats-this session'!`$A`$2:`$A`$5,'PLO Stats-this session'!`$K`$2:`$K`$5,2)"

# chart12.xml -> wsPLO ChartObject 5
Set-ChartSeriesFormula $wsPLO 5 1 "=SERIES('PLO Stats-this session'!R1,'PLO Stats-this session'!`$A`$2:`$A`$5,'PLO Stats-this session'!`$R`$2:`$R`$5,1)"

# chart13.xml -> wsPLO ChartObject 6
Set-ChartSeriesFormula $wsPLO 6 1 "=SERIES('PLO Stats-this session'!O1,'PLO Stats-this session'!`$A`$2:`$A`$5,'PLO Stats-this session'!`$O`$2:`$O`$5,1)"
Set-ChartSeriesFormula $wsPLO 6 2 "=SERIES('PLO Stats-this session'!P1,'PLO Stats-this session'!`$A`$2:`$A`$5,'PLO Stats-this session'!`$P`$2:`$P`$5,2)"

# chart14.xml -> wsPLO ChartObject 7
Set-ChartSeriesFormula $wsPLO 7 1 "=SERIES('PLO Stats-this session'!Q1,'PLO Stats-this session'!`$A`$2:`$A`$5,'PLO Stats-this session'!`$Q`$2:`$Q`$5,1)"
